$d = $word.ActiveDocument

# First paragraph: the hidden bookmark/id paragraph at top of doc.
$p1 = $d.Paragraphs(1)

# Add a paragraph border (5pt space on all sides) to match w:pBdr w:space="5".
$borders = $p1.Range.ParagraphFormat.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5

# Change the left indent from 120 twips (6pt) to 225 twips (11.25pt).
$p1.Format.LeftIndent = 11.25

# Replace the placeholder id text and swallow the trailing space run so the
# paragraph ends up with a single run and no trailing space.
$find = $p1.Range.Find
$find.Execute("**ID__AFFARS_pgi_5306_topic_3__ID** ", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_AF_PGI_5306_302_4__ID**", 2)
